# fix(publipostage): Add space before ":"
#
# 1) Every "statut_name" label in column B (rows 2-9) gets a space inserted
#    before the colon, e.g. "3: ..." -> "3 : ...".
# 2) Rows 4 and 5 (the NCT00978783 / NCT00320099 trial rows) swap places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix "X: ..." -> "X : ..." for every data row in column B ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 9 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $text = $cell.Value2
    if ($text -ne $null) {
        $newText = $text -replace '^(\d+):\s*', '$1 : '
        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}

# --- Step 2: swap the full contents of row 4 and row 5 ---
# (NCT00978783 / "Speech Effects..." row trades places with the
#  NCT00320099 / "Phase 3 Study of Corticotherapy..." row)
$cols = 1..9
$row4 = @{}
$row5 = @{}
foreach ($c in $cols) {
    $row4[$c] = $ws.Cells.Item(4, $c).Value2
    $row5[$c] = $ws.Cells.Item(5, $c).Value2
}
foreach ($c in $cols) {
    # Columns D/E are empty in both source rows - leave those cells alone
    # so no empty <c> node gets materialised where none existed before.
    if ($row4[$c] -eq $null -and $row5[$c] -eq $null) {
        continue
    }

    $cell4 = $ws.Cells.Item(4, $c)
    if ($row5[$c] -eq $null) {
        # Destination should end up blank - fully remove the cell rather
        # than writing an empty value, so it drops out like an untouched
        # empty cell would.
        $cell4.ClearContents()
    } else {
        # Every populated source cell here holds text, not a real number
        # (e.g. "3", "2009" are stored as strings), so briefly force the
        # cell to a text format before assigning the value back -
        # otherwise numeric-looking strings would silently turn into real
        # numbers. The style is restored to "Normal" right after so the
        # cell keeps its original (unstyled) appearance.
        $cell4.NumberFormat = "@"
        $cell4.Value = $row5[$c]
        $cell4.Style = "Normal"
    }

    $cell5 = $ws.Cells.Item(5, $c)
    if ($row4[$c] -eq $null) {
        $cell5.ClearContents()
    } else {
        $cell5.NumberFormat = "@"
        $cell5.Value = $row4[$c]
        $cell5.Style = "Normal"
    }
}
